# corrected data cleaning for pre/post/total fixation data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:R1): drop the bold/bordered/centered "header" style ---
# (font 1, border 1 and the cellXfs entry that used them go away; cells go
# back to the default style).
$ws.Range("A1:R1").ClearFormats()

# A1 no longer carries the "Unnamed: 0" label (pandas index column header).
$ws.Range("A1").ClearContents()

# --- Corrected numeric values -------------------------------------------
# Row 3 (Revisit count)
$ws.Range("C3").Value = 19
$ws.Range("E3").Value = 1
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 6
$ws.Range("I3").Value = 8
$ws.Range("N3").Value = 23
$ws.Range("O3").Value = 1

# Row 4 (Fixation count)
$ws.Range("C4").Value = 28
$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 8
$ws.Range("I4").Value = 10
$ws.Range("N4").Value = 72
$ws.Range("O4").Value = 2

# Row 5 (Dwell time (ms))
$ws.Range("C5").Value = 11235.07
$ws.Range("E5").Value = 495.83
$ws.Range("G5").Value = 1553
$ws.Range("H5").Value = 5001.89
$ws.Range("I5").Value = 4682.25
$ws.Range("N5").Value = 25567.58
$ws.Range("O5").Value = 924.98

# Row 6 (Dwell time (%))
$ws.Range("C6").Value = 15.04
$ws.Range("D6").Value = 1.06
$ws.Range("G6").Value = 2.08
$ws.Range("H6").Value = 6.7
$ws.Range("I6").Value = 6.27
$ws.Range("K6").Value = 0.53
$ws.Range("L6").Value = 5.92
$ws.Range("M6").Value = 0.18
$ws.Range("N6").Value = 34.24
$ws.Range("O6").Value = 1.24
$ws.Range("R6").Value = 0.18

# Row 7 (Fixation duration (ms))
$ws.Range("C7").Value = 401.25
$ws.Range("E7").Value = 247.91
$ws.Range("G7").Value = 310.6
$ws.Range("H7").Value = 625.24
$ws.Range("I7").Value = 468.23
$ws.Range("N7").Value = 355.11
$ws.Range("O7").Value = 462.49

# Row 8 (First fixation duration (ms))
$ws.Range("C8").Value = 132.16
$ws.Range("N8").Value = 74.22

# --- Drop the trailing empty rows 10-14 (dimension shrinks to A1:R9) ----
$ws.Rows("10:14").Delete()
